$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '329.44'
Set-TextValue "E2" '1.14%'
Set-TextValue "D3" '41.25'
Set-TextValue "E3" '4.14%'
Set-TextValue "D4" '5.645'
Set-TextValue "E4" '-0.41%'
Set-TextValue "D5" '0.08206'
Set-TextValue "E5" '2.17%'
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue "D6" '2.020'
Set-TextValue "E6" '-1.01%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue "D7" '8.757'
Set-TextValue "E7" '1.57%'
Set-TextValue "D8" '4.506'
Set-TextValue "E8" '0.43%'
Set-TextValue "D9" '2.965'
Set-TextValue "E9" '1.09%'
Set-TextValue "D10" '0.9226'
Set-TextValue "E10" '-0.02%'
Set-TextValue "E11" '2.63%'
Set-TextValue "D12" '0.1952'
Set-TextValue "E12" '-0.69%'
Set-TextValue "D13" '0.09334'
Set-TextValue "E13" '1.45%'
Set-TextValue "D14" '0.03840'
Set-TextValue "E14" '7.77%'
Set-TextValue "E15" '0.91%'
Set-TextValue "D16" '0.001310'
Set-TextValue "E16" '0.82%'
Set-TextValue "D17" '0.006151'
Set-TextValue "E17" '0.93%'
Set-TextValue "D19" '3.447'
Set-TextValue "E19" '2.86%'
Set-TextValue "E20" '-0.05%'
Set-TextValue "D21" '8.327'
Set-TextValue "E21" '-4.48%'
Set-TextValue "E22" '-0.35%'
Set-TextValue "E23" '6.28%'
Set-TextValue "D24" '0.04403'
Set-TextValue "E24" '0.38%'
Set-TextValue "E25" '-0.17%'
Set-TextValue "D26" '0.004317'
Set-TextValue "E26" '-6.22%'
Set-TextValue "D27" '0.0001201'
Set-TextValue "E27" '-2.40%'
Set-TextValue "D39" '0.02757'
Set-TextValue "E39" '9.62%'
Set-TextValue "D40" '0.05504'
Set-TextValue "E40" '3.20%'
Set-TextValue "D41" '0.007879'
Set-TextValue "E41" '5.34%'
Set-TextValue "D42" '0.1422'
Set-TextValue "E42" '1.19%'
Set-TextValue "D43" '0.008943'
Set-TextValue "E43" '-9.78%'
Set-TextValue "D44" '0.002171'
Set-TextValue "E44" '2.62%'
Set-TextValue "E45" '2.68%'
Set-TextValue "D46" '0.00006771'
Set-TextValue "E46" '1.37%'
Set-TextValue "E47" '0.00%'
Set-TextValue "D48" '0.003191'
Set-TextValue "E48" '7.24%'
Set-TextValue "E49" '0.00%'
Set-TextValue "D50" '0.00002101'
Set-TextValue "E50" '0.00%'
Set-TextValue "D51" '0.0002001'
Set-TextValue "E51" '0.00%'
